# Update "want to go" counts (column F) for the 展览 (Exhibition) sheet
# and the 全部类型 (All types) sheet, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# Map of row -> new F value for sheet "展览" (tab 1 / sheet1.xml)
$exhibitionUpdates = @{
    3  = 16480
    5  = 742
    6  = 15630
    7  = 74
    8  = 9292
    9  = 505
    12 = 132
    16 = 24
    17 = 97
    18 = 634
    20 = 18
    21 = 81
    22 = 1161
    25 = 35
    26 = 547
    32 = 70
    33 = 275
    34 = 378
    37 = 5756
}

# Map of row -> new F value for sheet "全部类型" (tab 4 / sheet4.xml)
# Rows 29/30 correspond to performance events not present in 展览,
# so rows after 26 are shifted by +2 relative to the 展览 sheet.
$allTypesUpdates = @{
    3  = 16480
    5  = 742
    6  = 15630
    7  = 74
    8  = 9292
    9  = 505
    12 = 132
    16 = 24
    17 = 97
    18 = 634
    20 = 18
    21 = 81
    22 = 1161
    25 = 35
    26 = 547
    34 = 70
    35 = 275
    36 = 378
    39 = 5756
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
